$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric to Excel but must stay stored as text
# (matches the source data, which keeps these as plain text strings).
$textCells = @(
    'D4',
    'D5',
    'D6',
    'D9',
    'D10',
    'D12',
    'D16',
    'D17',
    'D19',
    'D20',
    'D22',
    'D23',
    'D26',
    'D28',
    'D29',
    'D30',
    'D31',
    'D32',
    'D35',
    'D37',
    'D38',
    'D40',
    'D41',
    'D42',
    'D45',
    'D47',
    'D49',
    'D50',
    'D51'
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, in sheet order.
$ws.Range('D2').Value = '43.007.47'
$ws.Range('E2').Value = '  +0.91%  '
$ws.Range('D3').Value = '2.575.33'
$ws.Range('E3').Value = '  +2.26%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '315.01'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').Value = '100.40'
$ws.Range('E6').Value = '  +5.22%  '
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '0.539'
$ws.Range('E9').Value = '  +1.20%  '
$ws.Range('D10').Value = '36.28'
$ws.Range('E10').Value = '  +1.34%  '
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('D12').Value = '7.58'
$ws.Range('E12').Value = '  +0.64%  '
$ws.Range('D13').Value = '2.971.28'
$ws.Range('E13').Value = '  +2.19%  '
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').Value = '2.645.90'
$ws.Range('E15').Value = '  +5.78%  '
$ws.Range('D16').Value = '15.73'
$ws.Range('E16').Value = '  +3.13%  '
$ws.Range('D17').Value = '0.845'
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('D18').Value = '43.038.06'
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('D19').Value = '6.88'
$ws.Range('E19').Value = '  +2.44%  '
$ws.Range('D20').Value = '12.68'
$ws.Range('E20').Value = '  -1.12%  '
$ws.Range('D21').Value = '0.0₃0972'
$ws.Range('E21').Value = '  +1.36%  '
$ws.Range('D22').Value = '69.42'
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('D23').Value = '250.40'
$ws.Range('E23').Value = '  +0.26%  '
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('E25').Value = '  +0.80%  '
$ws.Range('D26').Value = '27.13'
$ws.Range('E26').Value = '  +2.17%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').Value = '2.41'
$ws.Range('E28').Value = '  -0.79%  '
$ws.Range('D29').Value = '40.67'
$ws.Range('E29').Value = '  -1.90%  '
$ws.Range('D30').Value = '10.33'
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('D31').Value = '5.86'
$ws.Range('E31').Value = '  -1.36%  '
$ws.Range('D32').Value = '157.58'
$ws.Range('E32').Value = '  -0.15%  '
$ws.Range('E33').Value = '  +4.38%  '
$ws.Range('E34').Value = '  +0.19%  '
$ws.Range('D35').Value = '0.0806'
$ws.Range('E35').Value = '  +3.41%  '
$ws.Range('E36').Value = '  -0.29%  '
$ws.Range('D37').Value = '18.86'
$ws.Range('E37').Value = '  -2.42%  '
$ws.Range('D38').Value = '2.53'
$ws.Range('E38').Value = '  +9.28%  '
$ws.Range('E39').Value = '  +1.16%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = '0.119'
$ws.Range('E40').Value = '  +0.63%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '23.85'
$ws.Range('E41').Value = '  +2.64%  '
$ws.Range('D42').Value = '4.05'
$ws.Range('E42').Value = '  +6.97%  '
$ws.Range('E43').Value = '  +0.24%  '
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').Value = '3.27'
$ws.Range('E45').Value = '  -1.72%  '
$ws.Range('D46').Value = '2.010.24'
$ws.Range('E46').Value = '  -1.39%  '
$ws.Range('D47').Value = '8.93'
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('D48').Value = '2.822.82'
$ws.Range('E48').Value = '  +2.10%  '
$ws.Range('D49').Value = '0.198'
$ws.Range('E49').Value = '  +3.10%  '
$ws.Range('D50').Value = '75.16'
$ws.Range('E50').Value = '  -0.02%  '
$ws.Range('D51').Value = '82.15'
$ws.Range('E51').Value = '  -2.57%  '
